$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (AC1) onto the
# three new header cells so they match the bold/centered/bordered style
# used by the rest of the header row (style index 1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Set the new header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 87
    $ws.Cells.Item($r, 31).Value = 75
    $ws.Cells.Item($r, 32).Value = 0
}
